# Update cryptocurrency price/volume data in the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these cells to Text format first so numeric-looking price strings
# (e.g. "1.000", "0.4490") are preserved exactly as text, not converted to numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = "30.471.15"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "2.106.10"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "332.42"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.5229"
$ws.Range("E7").Value = "  -1.82%  "
$ws.Range("D8").Value = "0.4490"
$ws.Range("E8").Value = "  +2.36%  "
$ws.Range("D9").Value = "53.78"
$ws.Range("E9").Value = "  +16.94%  "
$ws.Range("D10").Value = "0.08952"
$ws.Range("D11").Value = "1.160"
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("D12").Value = "24.46"
$ws.Range("E12").Value = "  -2.27%  "
$ws.Range("D13").Value = "2.099.44"
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("D14").Value = "6.759"
$ws.Range("D15").Value = "7.764"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").Value = "96.50"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "0.00001126"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").Value = "19.31"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "6.299"
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").Value = "30.526.87"
$ws.Range("E23").Value = "  -1.03%  "
$ws.Range("D24").Value = "12.33"
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").Value = "2.349"
$ws.Range("E25").Value = "  +4.35%  "
$ws.Range("D26").Value = "2.342.18"
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("E27").Value = "  -1.83%  "
$ws.Range("D28").Value = "2.590"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").Value = "163.61"
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("D31").Value = "1.201"
$ws.Range("E31").Value = "  +2.63%  "
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").Value = "1.680"
$ws.Range("E33").Value = "  +8.81%  "
$ws.Range("D34").Value = "6.162"
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("D35").Value = "3.939"
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("D36").Value = "10.41"
$ws.Range("E36").Value = "  +8.89%  "
$ws.Range("D37").Value = "0.02572"
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("D38").Value = "0.06788"
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("D39").Value = "12.87"
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").Value = "5.493"
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("D41").Value = "0.2281"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").Value = "0.6926"
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("D43").Value = "1.253"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "2.307"
$ws.Range("E45").Value = "  +3.35%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.6374"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "13.91"
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("D48").Value = "3.640"
$ws.Range("E48").Value = "  -0.87%  "
$ws.Range("D49").Value = "1.247"
$ws.Range("E49").Value = "  -2.37%  "
$ws.Range("E50").Value = "  +5.62%  "
$ws.Range("D51").Value = "83.16"
$ws.Range("E51").Value = "  +0.76%  "
